# cryptos.xlsx refresh -- GitHub Actions data pull (Sat Feb 10 22:31:23 UTC 2024).
# Coin/Price/Volume(1h) snapshot changed for the rows below. Rows 29-31 and
# 40-41 were additionally re-ranked, so Coin name + Link also move there.
#
# $forceText cells hold plain decimal-looking strings (e.g. "2.78", "0.139").
# Excel auto-converts such a string on a plain Range.Value assignment to a
# Number (and can even truncate a trailing zero, e.g. "247.60" -> 247.6), so
# those cells are pre-formatted as Text ("@") to keep them literal strings,
# matching every other Price cell already stored as text in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "47.714.62"; ForceText = $false },
    @{ Cell = "E2"; Value = "  +0.62%  "; ForceText = $false },

    @{ Cell = "D3"; Value = "2.494.26"; ForceText = $false },
    @{ Cell = "E3"; Value = "  -0.17%  "; ForceText = $false },

    @{ Cell = "E4"; Value = "  +0.05%  "; ForceText = $false },

    @{ Cell = "D5"; Value = "322.79"; ForceText = $true },
    @{ Cell = "E5"; Value = "  -0.40%  "; ForceText = $false },

    @{ Cell = "D6"; Value = "108.97"; ForceText = $true },
    @{ Cell = "E6"; Value = "  +0.93%  "; ForceText = $false },

    @{ Cell = "D7"; Value = "0.523"; ForceText = $true },
    @{ Cell = "E7"; Value = "  -0.79%  "; ForceText = $false },

    @{ Cell = "E8"; Value = "  +0.04%  "; ForceText = $false },

    @{ Cell = "D9"; Value = "0.552"; ForceText = $true },
    @{ Cell = "E9"; Value = "  +1.23%  "; ForceText = $false },

    @{ Cell = "D10"; Value = "40.51"; ForceText = $true },
    @{ Cell = "E10"; Value = "  +6.11%  "; ForceText = $false },

    @{ Cell = "D11"; Value = "0.0813"; ForceText = $true },
    @{ Cell = "E11"; Value = "  -0.16%  "; ForceText = $false },

    @{ Cell = "E12"; Value = "  +0.54%  "; ForceText = $false },

    @{ Cell = "D13"; Value = "18.68"; ForceText = $true },
    @{ Cell = "E13"; Value = "  +1.18%  "; ForceText = $false },

    @{ Cell = "D14"; Value = "7.21"; ForceText = $true },
    @{ Cell = "E14"; Value = "  -0.09%  "; ForceText = $false },

    @{ Cell = "D15"; Value = "2.890.66"; ForceText = $false },
    @{ Cell = "E15"; Value = "  +0.02%  "; ForceText = $false },

    @{ Cell = "D16"; Value = "2.492.99"; ForceText = $false },
    @{ Cell = "E16"; Value = "  -0.27%  "; ForceText = $false },

    @{ Cell = "D17"; Value = "0.852"; ForceText = $true },
    @{ Cell = "E17"; Value = "  -0.06%  "; ForceText = $false },

    @{ Cell = "D18"; Value = "47.652.49"; ForceText = $false },
    @{ Cell = "E18"; Value = "  +0.65%  "; ForceText = $false },

    @{ Cell = "D19"; Value = "13.21"; ForceText = $true },
    @{ Cell = "E19"; Value = "  +1.71%  "; ForceText = $false },

    @{ Cell = "D20"; Value = "6.63"; ForceText = $true },
    @{ Cell = "E20"; Value = "  -0.63%  "; ForceText = $false },

    @{ Cell = "D21"; Value = "2.78"; ForceText = $true },
    @{ Cell = "E21"; Value = "  +14.08%  "; ForceText = $false },

    @{ Cell = "D22"; Value = "0.0₃0942"; ForceText = $false },
    @{ Cell = "E22"; Value = "  -0.01%  "; ForceText = $false },

    @{ Cell = "D23"; Value = "70.82"; ForceText = $true },
    @{ Cell = "E23"; Value = "  +0.09%  "; ForceText = $false },

    @{ Cell = "D24"; Value = "247.60"; ForceText = $true },
    @{ Cell = "E24"; Value = "  -1.52%  "; ForceText = $false },

    @{ Cell = "D25"; Value = "2.56"; ForceText = $true },
    @{ Cell = "E25"; Value = "  -1.33%  "; ForceText = $false },

    @{ Cell = "E26"; Value = "  +0.11%  "; ForceText = $false },

    @{ Cell = "D27"; Value = "25.85"; ForceText = $true },
    @{ Cell = "E27"; Value = "  -1.60%  "; ForceText = $false },

    @{ Cell = "D28"; Value = "9.98"; ForceText = $true },
    @{ Cell = "E28"; Value = "  -1.07%  "; ForceText = $false },

    @{ Cell = "B29"; Value = "Kaspa"; ForceText = $false },
    @{ Cell = "C29"; Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; ForceText = $false },
    @{ Cell = "D29"; Value = "0.139"; ForceText = $true },
    @{ Cell = "E29"; Value = "  +1.26%  "; ForceText = $false },

    @{ Cell = "B30"; Value = "InjectiveProtocol"; ForceText = $false },
    @{ Cell = "C30"; Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; ForceText = $false },
    @{ Cell = "D30"; Value = "35.15"; ForceText = $true },
    @{ Cell = "E30"; Value = "  -0.67%  "; ForceText = $false },

    @{ Cell = "B31"; Value = "Toncoin"; ForceText = $false },
    @{ Cell = "C31"; Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; ForceText = $false },
    @{ Cell = "D31"; Value = "2.08"; ForceText = $true },
    @{ Cell = "E31"; Value = "  -0.32%  "; ForceText = $false },

    @{ Cell = "D32"; Value = "49.75"; ForceText = $true },
    @{ Cell = "E32"; Value = "  +0.57%  "; ForceText = $false },

    @{ Cell = "D33"; Value = "20.01"; ForceText = $true },
    @{ Cell = "E33"; Value = "  +1.00%  "; ForceText = $false },

    @{ Cell = "D34"; Value = "5.36"; ForceText = $true },
    @{ Cell = "E34"; Value = "  -3.20%  "; ForceText = $false },

    @{ Cell = "D35"; Value = "0.0791"; ForceText = $true },
    @{ Cell = "E35"; Value = "  -0.48%  "; ForceText = $false },

    @{ Cell = "E36"; Value = "  +0.16%  "; ForceText = $false },

    @{ Cell = "D37"; Value = "1.96"; ForceText = $true },
    @{ Cell = "E37"; Value = "  -1.93%  "; ForceText = $false },

    @{ Cell = "E38"; Value = "  -0.83%  "; ForceText = $false },

    @{ Cell = "D39"; Value = "2.96"; ForceText = $true },
    @{ Cell = "E39"; Value = "  -1.49%  "; ForceText = $false },

    @{ Cell = "B40"; Value = "EnergySwap"; ForceText = $false },
    @{ Cell = "C40"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; ForceText = $false },
    @{ Cell = "D40"; Value = "22.61"; ForceText = $true },
    @{ Cell = "E40"; Value = "  +7.39%  "; ForceText = $false },

    @{ Cell = "B41"; Value = "Stellar"; ForceText = $false },
    @{ Cell = "C41"; Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; ForceText = $false },
    @{ Cell = "D41"; Value = "0.112"; ForceText = $true },
    @{ Cell = "E41"; Value = "  -0.28%  "; ForceText = $false },

    @{ Cell = "E42"; Value = "  -1.00%  "; ForceText = $false },

    @{ Cell = "D43"; Value = "119.28"; ForceText = $true },
    @{ Cell = "E43"; Value = "  -1.99%  "; ForceText = $false },

    @{ Cell = "E44"; Value = "  +0.11%  "; ForceText = $false },

    @{ Cell = "D45"; Value = "2.002.13"; ForceText = $false },
    @{ Cell = "E45"; Value = "  +1.63%  "; ForceText = $false },

    @{ Cell = "D46"; Value = "3.06"; ForceText = $true },
    @{ Cell = "E46"; Value = "  +0.79%  "; ForceText = $false },

    @{ Cell = "E47"; Value = "  -3.45%  "; ForceText = $false },

    @{ Cell = "D48"; Value = "1.81"; ForceText = $true },
    @{ Cell = "E48"; Value = "  +0.23%  "; ForceText = $false },

    @{ Cell = "D49"; Value = "9.03"; ForceText = $true },
    @{ Cell = "E49"; Value = "  -0.45%  "; ForceText = $false },

    @{ Cell = "D50"; Value = "5.17"; ForceText = $true },
    @{ Cell = "E50"; Value = "  -2.40%  "; ForceText = $false },

    @{ Cell = "D51"; Value = "56.86"; ForceText = $true },
    @{ Cell = "E51"; Value = "  +2.94%  "; ForceText = $false }
)

foreach ($u in $updates) {
    if ($u.ForceText) {
        # Lock the cell to Text format *before* assigning, so the numeric-
        # looking string is stored verbatim instead of being parsed as a Number.
        $ws.Range($u.Cell).NumberFormat = "@"
    }
    $ws.Range($u.Cell).Value = $u.Value
}
